$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Done" column for the newly-completed problems (rows 6-13 and 15-19)
$rows = @(6,7,8,9,10,11,12,13,15,16,17,18,19)
foreach ($r in $rows) {
    $ws.Range("C$r").Value = "<yes>"
}

# Update the active selection to match the author's final cursor position
$ws.Range("D18").Select()
